# Applies the "added point parameter example and USGS library parameter
# values for reference" edit:
#   - inserts a new "BIO" browse-product row (between HEM and HY2)
#   - inserts a new "MAF2" browse-product row (between MAF and PAL)
#   - updates the sheet dimension / filter-database defined name to match
#     the new row count
#   - updates the view (zoom + selection) to match the saved state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the "BIO" row at row 6 (pushes HY2..TRU down by one) -----------
$ws.Rows.Item(6).Insert()
$ws.Range("A6").Value = "BIO"
$ws.Range("B6").Value = "BD1200"
$ws.Range("C6").Value = "BD670"
$ws.Range("D6").Value = "D700"

# --- Insert the "MAF2" row at row 12 (pushes PAL..TRU down by one) ---------
$ws.Rows.Item(12).Insert()
$ws.Range("A12").Value = "MAF2"
$ws.Range("B12").Value = "OLINDEX3"
$ws.Range("C12").Value = "HCPINDEX2"
$ws.Range("D12").Value = "BD920_2"

# --- Defined name _FilterDatabase now spans the two extra rows -------------
$fd = $wb.Names.Item(1)
$fd.RefersTo = "=Sheet1!`$A`$1:`$D`$17"

# --- View state: zoom + active selection ------------------------------------
$excel.ActiveWindow.Zoom = 184
[void]$ws.Range("B7").Select()
